$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.024.64"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "3.989.20"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "595.23"
$ws.Range("E5").Value = "  +10.00%  "
$ws.Range("D6").Value = "161.50"
$ws.Range("E6").Value = "  +7.75%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "54.13"
$ws.Range("E11").Value = "  -4.86%  "
$ws.Range("D12").Value = "0.0000320"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "10.97"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "4.620.71"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.995.89"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +8.03%  "
$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "20.38"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "72.627.03"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").Value = "434.11"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  +13.21%  "
$ws.Range("D23").Value = "96.25"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "4.38"
$ws.Range("E26").Value = "  +15.97%  "
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "36.41"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").Value = "13.74"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "48.88"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").Value = "671.76"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").Value = "70.68"
$ws.Range("E36").Value = "  +7.67%  "
$ws.Range("D37").Value = "0.0₃0893"
$ws.Range("E37").Value = "  +8.73%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "0.0492"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  +9.38%  "
$ws.Range("D46").Value = "0.149"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "3.44"
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").Value = "2.860.13"
$ws.Range("E49").Value = "  +4.68%  "
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("E51").Value = "  +3.88%  "
